$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.736.89'
$ws.Range('E2').Value = '  +0.86%  '
$ws.Range('D3').Value = '1.645.66'
$ws.Range('E3').Value = '  +1.15%  '
$ws.Range('E4').Value = '  +0.31%  '
$ws.Range('D5').Formula = '="216.27"'
$ws.Range('D5').Copy()
$ws.Range('D5').PasteSpecial(-4163)
$ws.Range('E6').Value = '  +0.60%  '
$ws.Range('E7').Value = '  +0.25%  '
$ws.Range('E8').Value = '  +0.81%  '
$ws.Range('E9').Value = '  +0.46%  '
$ws.Range('D10').Formula = '="19.21"'
$ws.Range('D10').Copy()
$ws.Range('D10').PasteSpecial(-4163)
$ws.Range('E10').Value = '  +2.18%  '
$ws.Range('D11').Formula = '="0.0842"'
$ws.Range('D11').Copy()
$ws.Range('D11').PasteSpecial(-4163)
$ws.Range('E11').Value = '  -0.23%  '
$ws.Range('D12').Value = '1.874.53'
$ws.Range('E12').Value = '  +1.14%  '
$ws.Range('D13').Value = '1.636.86'
$ws.Range('E13').Value = '  -0.63%  '
$ws.Range('E14').Value = '  +1.09%  '
$ws.Range('D15').Formula = '="0.532"'
$ws.Range('D15').Copy()
$ws.Range('D15').PasteSpecial(-4163)
$ws.Range('E15').Value = '  +1.64%  '
$ws.Range('D16').Formula = '="65.28"'
$ws.Range('D16').Copy()
$ws.Range('D16').PasteSpecial(-4163)
$ws.Range('E16').Value = '  +0.25%  '
$ws.Range('D17').Value = '26.749.85'
$ws.Range('E17').Value = '  +0.79%  '
$ws.Range('D18').Value = '0.0₃0743'
$ws.Range('E18').Value = '  +0.38%  '
$ws.Range('D19').Formula = '="218.11"'
$ws.Range('D19').Copy()
$ws.Range('D19').PasteSpecial(-4163)
$ws.Range('E19').Value = '  +1.68%  '
$ws.Range('E20').Value = '  +0.26%  '
$ws.Range('D21').Formula = '="4.36"'
$ws.Range('D21').Copy()
$ws.Range('D21').PasteSpecial(-4163)
$ws.Range('E21').Value = '  +1.76%  '
$ws.Range('D22').Formula = '="2.44"'
$ws.Range('D22').Copy()
$ws.Range('D22').PasteSpecial(-4163)
$ws.Range('E22').Value = '  +17.80%  '
$ws.Range('E23').Value = '  -0.10%  '
$ws.Range('E24').Value = '  +1.85%  '
$ws.Range('D25').Formula = '="146.41"'
$ws.Range('D25').Copy()
$ws.Range('D25').PasteSpecial(-4163)
$ws.Range('E25').Value = '  -0.91%  '
$ws.Range('E26').Value = '  +0.30%  '
$ws.Range('E27').Value = '  -0.05%  '
$ws.Range('E28').Value = '  +3.87%  '
$ws.Range('D29').Formula = '="15.74"'
$ws.Range('D29').Copy()
$ws.Range('D29').PasteSpecial(-4163)
$ws.Range('E29').Value = '  +1.44%  '
$ws.Range('E30').Value = '  +1.64%  '
$ws.Range('E31').Value = '  +1.42%  '
$ws.Range('E32').Value = '  -0.14%  '
$ws.Range('D33').Formula = '="3.00"'
$ws.Range('D33').Copy()
$ws.Range('D33').PasteSpecial(-4163)
$ws.Range('E33').Value = '  +1.31%  '
$ws.Range('D34').Value = '1.277.63'
$ws.Range('E34').Value = '  +2.55%  '
$ws.Range('E35').Value = '  +2.96%  '
$ws.Range('E36').Value = '  +2.97%  '
$ws.Range('E37').Value = '  +1.92%  '
$ws.Range('D38').Formula = '="0.538"'
$ws.Range('D38').Copy()
$ws.Range('D38').PasteSpecial(-4163)
$ws.Range('E38').Value = '  +5.74%  '
$ws.Range('D39').Formula = '="0.829"'
$ws.Range('D39').Copy()
$ws.Range('D39').PasteSpecial(-4163)
$ws.Range('E39').Value = '  +4.44%  '
$ws.Range('E40').Value = '  +0.29%  '
$ws.Range('D41').Formula = '="0.815"'
$ws.Range('D41').Copy()
$ws.Range('D41').PasteSpecial(-4163)
$ws.Range('E41').Value = '  +1.93%  '
$ws.Range('D42').Formula = '="2.25"'
$ws.Range('D42').Copy()
$ws.Range('D42').PasteSpecial(-4163)
$ws.Range('E42').Value = '  -1.30%  '
$ws.Range('D43').Formula = '="5.45"'
$ws.Range('D43').Copy()
$ws.Range('D43').PasteSpecial(-4163)
$ws.Range('E43').Value = '  +2.14%  '
$ws.Range('D44').Value = '1.786.02'
$ws.Range('E44').Value = '  +1.24%  '
$ws.Range('D45').Formula = '="91.99"'
$ws.Range('D45').Copy()
$ws.Range('D45').PasteSpecial(-4163)
$ws.Range('E45').Value = '  -1.41%  '
$ws.Range('D46').Formula = '="59.69"'
$ws.Range('D46').Copy()
$ws.Range('D46').PasteSpecial(-4163)
$ws.Range('E46').Value = '  +8.82%  '
$ws.Range('B48').Value = 'Cronos'
$ws.Range('C48').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D48').Formula = '="0.0515"'
$ws.Range('D48').Copy()
$ws.Range('D48').PasteSpecial(-4163)
$ws.Range('E48').Value = '  +1.04%  '
$ws.Range('B49').Value = 'BabyDogeCoin'
$ws.Range('C49').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D49').Value = '0.0₇0996'
$ws.Range('E49').Value = '  -4.16%  '
$ws.Range('D50').Formula = '="7.74"'
$ws.Range('D50').Copy()
$ws.Range('D50').PasteSpecial(-4163)
$ws.Range('E50').Value = '  +3.18%  '
$ws.Range('D51').Formula = '="0.0971"'
$ws.Range('D51').Copy()
$ws.Range('D51').PasteSpecial(-4163)
$ws.Range('E51').Value = '  +1.21%  '
